$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

# New rows to append (Date, Timestamp, Hour, Location, Value, Status)
$rows = @(
    @("2026-02-01", "16:01:02", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:01:10", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:01:21", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:01:31", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:01:42", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:01:52", "16:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 55

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]

    # Column A holds a date-shaped string (e.g. "2026-02-01"). Assigning it directly
    # as .Value causes Excel to auto-convert it into a real date serial number with a
    # new number-format style, which does not match the source data (plain text).
    # Entering it as a text formula and then converting the formula to its plain
    # computed value keeps it as literal text with no special formatting/style.
    $ws.Cells.Item($r, 1).Formula = '="' + $rowData[0] + '"'

    # Remaining columns are plain text already and can be assigned directly.
    for ($c = 2; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

$dateRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $rows.Count - 1, 1))
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false
